$d = $word.ActiveDocument

# The document has an existing blank paragraph (4th body paragraph) right
# after the "startup.cs file." paragraph and before the "How to use the E
# Shopping Cart" heading. We want to insert three new list paragraphs
# between the "startup.cs file." paragraph and that existing blank
# paragraph, without disturbing the blank paragraph itself.
#
# InsertXML on a zero-length (collapsed) range is ambiguous about which
# neighbouring paragraph "owns" the insertion point, and can end up
# replacing the wrong paragraph. So instead we first call
# InsertParagraphBefore() to create a brand-new blank paragraph right
# before the existing blank one (this leaves the existing blank paragraph,
# and everything after it, completely untouched), then we target that
# freshly created paragraph's whole range with InsertXML to fill it with
# our three new list paragraphs.
$anchorPara = $d.Paragraphs(4)
$anchorRange = $anchorPara.Range
$anchorRange.Collapse(1)
$anchorRange.InsertParagraphBefore()

# Paragraph 4 is now the freshly created (still empty) paragraph; the
# original blank paragraph has shifted down to become paragraph 5.
$insertRange = $d.Paragraphs(4).Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Created Models folder and in it created two model classes &#8211; Item.cs and Category.cs. </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>In order to process the data in models</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>,</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t xml:space="preserve"> repositories are needed. </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>Since DI is used</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>,</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t xml:space="preserve"> Interfaces are created for each repository first. Then register the interfaces with the corresponding repository as a service and inject them into the application.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Created </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Interfaces &#8211; created </w:t>
  </w:r>
  <w:r>
    <w:t>Item Repository</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> Interface(IItemRepository)</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>and Category</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> Repository</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> Interface(I</w:t>
  </w:r>
  <w:r>
    <w:t>Category</w:t>
  </w:r>
  <w:r>
    <w:t>Repository)</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>in Models folder</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> with the relevant methods in them.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
</w:p>
'@

$insertRange.InsertXML($xml)
